$d = $word.ActiveDocument

$replacements = @(
    @("535÷6=89, 1", "514÷6=85, 4"),
    @("387÷4=96, 3", "474÷3=158, 0"),
    @("406÷9=45, 1", "666÷2=333, 0"),
    @("370÷7=52, 6", "386÷3=128, 2"),
    @("895÷9=99, 4", "671÷6=111, 5"),
    @("850÷4=212, 2", "769÷4=192, 1"),
    @("376÷4=94, 0", "749÷7=107, 0"),
    @("546÷2=273, 0", "738÷2=369, 0"),
    @("568÷8=71, 0", "459÷9=51, 0"),
    @("385÷5=77, 0", "522÷3=174, 0"),
    @("537÷5=107, 2", "821÷9=91, 2"),
    @("791÷4=197, 3", "415÷3=138, 1"),
    @("521÷4=130, 1", "483÷4=120, 3"),
    @("161÷5=32, 1", "269÷7=38, 3"),
    @("698÷2=349, 0", "298÷7=42, 4"),
    @("103÷8=12, 7", "540÷6=90, 0"),
    @("771÷5=154, 1", "523÷9=58, 1"),
    @("441÷8=55, 1", "629÷7=89, 6"),
    @("180÷7=25, 5", "782÷2=391, 0"),
    @("480÷2=240, 0", "235÷8=29, 3"),
    @("438÷5=87, 3", "757÷9=84, 1"),
    @("272÷5=54, 2", "715÷5=143, 0"),
    @("973÷5=194, 3", "245÷7=35, 0"),
    @("445÷8=55, 5", "373÷4=93, 1"),
    @("214÷9=23, 7", "857÷3=285, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
